# Agregados los datos de semana de lunes 22 de enero
# Adds a new "28_01_2024" column (F) with the days-worked counts for each
# recepcionista, mirroring the existing weekly columns (B..E).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header for the week of 28_01_2024
$ws.Range("F1").Value = "28_01_2024"

# New weekly data per recepcionista (rows 2-5 -> Alejandro, Camila, Betty, Felipe)
$ws.Range("F2").Value = 5
$ws.Range("F3").Value = 3
$ws.Range("F4").Value = 3
$ws.Range("F5").Value = 4

# Reflect the new active selection after entering the data
$ws.Range("F3").Select() | Out-Null
